$d = $word.ActiveDocument

# --- Step 1: fix the CREATE DATABASE statement text ---------------------
# " CHARACTER SET utf8 COLLATE DEFAULT;" -> " CHARACTER SET utf8 ;"
$findRng = $d.Content
$null = $findRng.Find.Execute(" CHARACTER SET utf8 COLLATE DEFAULT;", $true, $false, $false, $false, $false, $true, 1, $false, " CHARACTER SET utf8 ;", 2)

# --- Step 2: relocate the "_GoBack" bookmark -----------------------------
# In the original document the bookmark sits (collapsed) in the empty
# paragraph right after the final ");" of the script. Word actually wants
# it to span from the start of the "CREATE DATABASE ..." paragraph through
# the end of that final ");" paragraph (this is simply where the last edit
# in the document happened to occur).

# Start boundary: beginning of the paragraph that starts the SQL script.
$startRng = $d.Content
$null = $startRng.Find.Execute("CREATE DATABASE")
$startPos = $startRng.Start

# End boundary: end of the paragraph that sits right before the paragraph
# which currently holds the "_GoBack" bookmark (i.e. the paragraph with the
# closing ");").
$bmRange = $d.Bookmarks("_GoBack").Range
$paraCount = $d.Paragraphs.Count
$bmParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $bmRange.Start -and $bmRange.Start -lt $p.Range.End) {
        $bmParaIndex = $i
        break
    }
}
$endPara = $d.Paragraphs.Item($bmParaIndex - 1)
$endPos = $endPara.Range.End - 1

# Move the bookmark: delete the old one, add a new one spanning the range.
$d.Bookmarks("_GoBack").Delete()
$newTarget = $d.Range($startPos, $endPos)
$d.Bookmarks.Add("_GoBack", $newTarget)
